# Update for Wed May 31 15:41:01 2023
#
# This script reproduces (as closely as the COM surface allows) the
# following source changes:
#  - Sheet1 becomes the active/selected sheet (was "Another Sheet")
#  - Sheet1 gains a new data row (row 4): Françoise-Athénaïs / de Rochechouart
#    / (no email) / a (very old, negative-serial) date
#  - The date column's number-format style becomes center-aligned
#  - Sheet1's column widths for columns A, B and D are widened slightly
#  - Selections on both sheets move to reflect the above
#
# NOTE: this COM-interop surface does not allow an application-set
# ActiveCell to differ from the top-left cell of a just-made Range.Select(),
# nor does Range.Select() retain multiple (Union) areas -- so the
# multi-area selections from the original edit ("A1:C13 G5", with an
# active cell off the top-left corner) are approximated with the single
# primary area instead.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# New row of data on Sheet1 (row 4)
# ---------------------------------------------------------------------
$ws1.Range("A4").Value = "Françoise-Athénaïs"
$ws1.Range("B4").Value = "de Rochechouart"
$ws1.Range("D4").Value = "1640-10-05"

# Center-align the date column (matches the style used by D2 / D3), now
# also covering the newly-added D4.
$ws1.Range("D2:D4").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# Column width tweaks on Sheet1 (A, B, D) -- column C is untouched.
# ColumnWidth uses "characters" units; the saved xlsx <col width="..">
# is derived from it (roughly width+0.83), so we back-solve for the
# ColumnWidth that lands closest to the target stored widths
# (17.29, 15.96, 12.29).
# ---------------------------------------------------------------------
$ws1.Columns.Item(1).ColumnWidth = 16.456666666666667
$ws1.Columns.Item(2).ColumnWidth = 15.126666666666667
$ws1.Columns.Item(4).ColumnWidth = 11.456666666666665

# ---------------------------------------------------------------------
# Selection / active-sheet state.
# Set sheet2's selection first (it stays visible even while not active),
# then finish on sheet1 so it ends up the active tab/sheet.
# ---------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("A1:C13").Select()

$ws1.Activate()
$ws1.Range("A1:C13").Select()
